$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 "Uebungen" ---
# Row 2: Push Ups
$ws1.Range("A2").Value = "Push Ups"
$ws1.Range("B2").Value = 3
$ws1.Range("C2").Value = 10

# Row 3: Pull Ups
$ws1.Range("A3").Value = "Pull Ups"
$ws1.Range("B3").Value = 3
$ws1.Range("C3").Value = 8

# Row 4: Squats (new)
$ws1.Range("A4").Value = "Squats"
$ws1.Range("B4").Value = 3
$ws1.Range("C4").Value = 10

# Row 5: Deadlifts (new)
$ws1.Range("A5").Value = "Deadlifts"
$ws1.Range("B5").Value = 3
$ws1.Range("C5").Value = 10

# --- Sheet2 "Trainingsstatistiken" ---
$ws2.Range("A1").Value = "Datum"
$ws2.Range("B1").Value = "Push Ups"
$ws2.Range("C1").Value = "Pull Ups"
$ws2.Range("D1").Value = "Squats"
$ws2.Range("E1").Value = "Deadlifts"

# --- Header formatting: thin bottom border on header rows ---
$ws1.Range("A1:C1").Borders.Item(9).LineStyle = 1
$ws2.Range("A1").Borders.Item(9).LineStyle = 1
